$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") that changed.
$updates = @{
    2  = 8022
    3  = 7644
    4  = 112
    8  = 127
    10 = 153
    12 = 691
    13 = 114
    14 = 1221
    16 = 44
    19 = 103
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
